$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A fresh scrape was appended: insert a new column right before the "nom"
# column (currently W) so "nom" shifts to X and "url_produit" shifts to Y.
$ws.Columns("W:W").Insert()

# Header for the newly-inserted price/timestamp column.
$ws.Range("W1").Value = "2026-01-28 16:22:56"

# For every product row that already had a price history, carry the most
# recent known price (previously the last column, V) forward into the new
# column, since the price did not change between scrapes.
$ws.Range("V2:V80").Copy()
$ws.Range("W2:W80").PasteSpecial(-4163)
$excel.CutCopyMode = 0
